# Handback status report regeneration:
# - the "06605d30-..." handback replaces the previous "cb1d4349-..." run for
#   the already-tracked source file (new xliff hashes / timestamps), and
# - a brand-new source file "32385396-..." is handed back for the first
#   time, so a second data row gets appended to every sheet/table.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # OLE (BGR) encoding of RGB 6495ED - the workbook's existing HyperLink font color

function Set-HyperlinkLook($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

function Set-DateLook($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Hyperlinks get rebuilt from scratch (clears the one stale link so it can be
# re-created pointing at the refreshed file, then adds the new row's link).
$ov.Range("A1").Hyperlinks.Delete() | Out-Null

# Row 2 (existing tracked file) gets refreshed metadata
$ov.Range("A2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.md"
$ov.Range("B2").Value = "e2e\06605d30-a0a8-472b-8fa3-6f10e3bae56c.md"
$ov.Range("G2").Value = "2017-02-09 09:38:49"
Set-DateLook $ov.Range("G2")

$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de0e3e82db7bcddd61cb0e68ce1d15131af63e/e2e/06605d30-a0a8-472b-8fa3-6f10e3bae56c.md", [Type]::Missing, [Type]::Missing, "e2e\06605d30-a0a8-472b-8fa3-6f10e3bae56c.md") | Out-Null
Set-HyperlinkLook $ov.Range("B2")

# Row 3: brand new handback entry
$ov.Range("A3").Value = "32385396-5055-45ef-8731-9774dff4db0e.md"
$ov.Range("B3").Value = "e2e\32385396-5055-45ef-8731-9774dff4db0e.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = "2017-02-09 09:37:35"
Set-DateLook $ov.Range("G3")

$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de0e3e82db7bcddd61cb0e68ce1d15131af63e/e2e/32385396-5055-45ef-8731-9774dff4db0e.md", [Type]::Missing, [Type]::Missing, "e2e\32385396-5055-45ef-8731-9774dff4db0e.md") | Out-Null
Set-HyperlinkLook $ov.Range("B3")

$ov.ListObjects.Item(1).Resize($ov.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A1").Hyperlinks.Delete() | Out-Null

# Row 2 refresh
$zh.Range("A2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.md"
$zh.Range("G2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.230e561e228ff076c6c011518a567862c8c7c783.zh-cn.xlf"
$zh.Range("H2").Value = "2017-02-09 09:38:29"
Set-DateLook $zh.Range("H2")
$zh.Range("J2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.md"
$zh.Range("K2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.230e561e228ff076c6c011518a567862c8c7c783.zh-cn.xlf"
$zh.Range("L2").Value = "2017-02-09 09:39:38"
Set-DateLook $zh.Range("L2")

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de0e3e82db7bcddd61cb0e68ce1d15131af63e/e2e/06605d30-a0a8-472b-8fa3-6f10e3bae56c.md", [Type]::Missing, [Type]::Missing, "06605d30-a0a8-472b-8fa3-6f10e3bae56c.md") | Out-Null
Set-HyperlinkLook $zh.Range("A2")
$zh.Hyperlinks.Add($zh.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/acd7bc07c9f275dd436eb0dd990111c68b055e74/e2e/06605d30-a0a8-472b-8fa3-6f10e3bae56c.md", [Type]::Missing, [Type]::Missing, "06605d30-a0a8-472b-8fa3-6f10e3bae56c.md") | Out-Null
Set-HyperlinkLook $zh.Range("J2")

# Row 3: brand new handback entry
$zh.Range("A3").Value = "32385396-5055-45ef-8731-9774dff4db0e.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "True"
$zh.Range("G3").Value = "32385396-5055-45ef-8731-9774dff4db0e.419f41cdd1dea672225752af5f50b10dc1def735.zh-cn.xlf"
$zh.Range("H3").Value = "2017-02-09 09:37:11"
Set-DateLook $zh.Range("H3")
$zh.Range("J3").Value = "32385396-5055-45ef-8731-9774dff4db0e.md"
$zh.Range("K3").Value = "32385396-5055-45ef-8731-9774dff4db0e.419f41cdd1dea672225752af5f50b10dc1def735.zh-cn.xlf"
$zh.Range("L3").Value = "2017-02-09 09:39:38"
Set-DateLook $zh.Range("L3")
$zh.Range("O3").Value = "True"
$zh.Range("Q3").Value = "False"

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de0e3e82db7bcddd61cb0e68ce1d15131af63e/e2e/32385396-5055-45ef-8731-9774dff4db0e.md", [Type]::Missing, [Type]::Missing, "32385396-5055-45ef-8731-9774dff4db0e.md") | Out-Null
Set-HyperlinkLook $zh.Range("A3")
$zh.Hyperlinks.Add($zh.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/acd7bc07c9f275dd436eb0dd990111c68b055e74/e2e/32385396-5055-45ef-8731-9774dff4db0e.md", [Type]::Missing, [Type]::Missing, "32385396-5055-45ef-8731-9774dff4db0e.md") | Out-Null
Set-HyperlinkLook $zh.Range("J3")

$zh.ListObjects.Item(1).Resize($zh.Range("A1:R3"))

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A1").Hyperlinks.Delete() | Out-Null

# Row 2 refresh
$de.Range("A2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.md"
$de.Range("G2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.230e561e228ff076c6c011518a567862c8c7c783.de-de.xlf"
$de.Range("L2").Value = "2017-02-09 09:40:07"
Set-DateLook $de.Range("L2")
$de.Range("J2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.md"
$de.Range("K2").Value = "06605d30-a0a8-472b-8fa3-6f10e3bae56c.230e561e228ff076c6c011518a567862c8c7c783.de-de.xlf"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de0e3e82db7bcddd61cb0e68ce1d15131af63e/e2e/06605d30-a0a8-472b-8fa3-6f10e3bae56c.md", [Type]::Missing, [Type]::Missing, "06605d30-a0a8-472b-8fa3-6f10e3bae56c.md") | Out-Null
Set-HyperlinkLook $de.Range("A2")
$de.Hyperlinks.Add($de.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b645be74980d28c81422f403ce864bacd3e459a0/e2e/06605d30-a0a8-472b-8fa3-6f10e3bae56c.md", [Type]::Missing, [Type]::Missing, "06605d30-a0a8-472b-8fa3-6f10e3bae56c.md") | Out-Null
Set-HyperlinkLook $de.Range("J2")

# Row 3: brand new handback entry
$de.Range("A3").Value = "32385396-5055-45ef-8731-9774dff4db0e.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "True"
$de.Range("G3").Value = "32385396-5055-45ef-8731-9774dff4db0e.419f41cdd1dea672225752af5f50b10dc1def735.de-de.xlf"
$de.Range("H3").Value = "2017-02-09 09:37:35"
Set-DateLook $de.Range("H3")
$de.Range("J3").Value = "32385396-5055-45ef-8731-9774dff4db0e.md"
$de.Range("K3").Value = "32385396-5055-45ef-8731-9774dff4db0e.419f41cdd1dea672225752af5f50b10dc1def735.de-de.xlf"
$de.Range("L3").Value = "2017-02-09 09:40:07"
Set-DateLook $de.Range("L3")
$de.Range("O3").Value = "True"
$de.Range("Q3").Value = "False"

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de0e3e82db7bcddd61cb0e68ce1d15131af63e/e2e/32385396-5055-45ef-8731-9774dff4db0e.md", [Type]::Missing, [Type]::Missing, "32385396-5055-45ef-8731-9774dff4db0e.md") | Out-Null
Set-HyperlinkLook $de.Range("A3")
$de.Hyperlinks.Add($de.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b645be74980d28c81422f403ce864bacd3e459a0/e2e/32385396-5055-45ef-8731-9774dff4db0e.md", [Type]::Missing, [Type]::Missing, "32385396-5055-45ef-8731-9774dff4db0e.md") | Out-Null
Set-HyperlinkLook $de.Range("J3")

$de.ListObjects.Item(1).Resize($de.Range("A1:R3"))
